$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 112; existing rows 112+ shift down to 114+
$ws.Rows("112:113").Insert()

# Row 112: new Ajo / Chino / Primera entry dated 44634 ($/caja 10 kilos)
$ws.Cells.Item(112, 1).Value = 9
$ws.Cells.Item(112, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(112, 3).Value = "Metropolitana"
$ws.Cells.Item(112, 4).Value = 44634
$ws.Cells.Item(112, 5).Value = 13
$ws.Cells.Item(112, 6).Value = 100112003
$ws.Cells.Item(112, 7).Value = "Ajo"
$ws.Cells.Item(112, 8).Value = "Chino"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 520
$ws.Cells.Item(112, 11).Value = 17500
$ws.Cells.Item(112, 12).Value = 18000
$ws.Cells.Item(112, 13).Value = 17750
$ws.Cells.Item(112, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(112, 15).Value = "China"
$ws.Cells.Item(112, 16).Value = 1775
$ws.Cells.Item(112, 17).Value = 10
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Row 113: new Ajo / Chino / Primera entry dated 44634 ($/malla 10 kilos)
$ws.Cells.Item(113, 1).Value = 9
$ws.Cells.Item(113, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(113, 3).Value = "Metropolitana"
$ws.Cells.Item(113, 4).Value = 44634
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = 100112003
$ws.Cells.Item(113, 7).Value = "Ajo"
$ws.Cells.Item(113, 8).Value = "Chino"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 250
$ws.Cells.Item(113, 11).Value = 18000
$ws.Cells.Item(113, 12).Value = 18500
$ws.Cells.Item(113, 13).Value = 18250
$ws.Cells.Item(113, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(113, 15).Value = "China"
$ws.Cells.Item(113, 16).Value = 1825
$ws.Cells.Item(113, 17).Value = 10
$ws.Cells.Item(113, 18).Value = "Hortaliza"
